$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("B18").Value = "sed_calcCorg =(1-dum_por) * loc_fPOC / loc_new_sed_vol"
$ws.Range("F18").Value = "!!!! Have to multiply (1-por), NOT devide!!!"
$ws.Range("C6").Select() | Out-Null
